$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with default (unstyled) formatting, never modified by this script,
# used to restore original cell formatting after forcing text values below.
$fmtSource = $ws.Range("B4")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.934.00'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.33%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.410.48'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.83'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -1.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.03'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.39%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.410.99'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.475'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.56'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.80%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.992.56'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.34'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000175'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.406.73'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.22%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.051.45'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.31'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.59%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.38'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.31'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '389.17'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.567'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '72.85'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.551.56'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.179'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.60%  '
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.19%  '
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.41'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -4.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.19'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.25%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.18'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.90'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.02'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.438.26'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.12'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -2.26%  '
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '167.58'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.89%  '
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.55'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0786'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.30'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.42%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.793'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.50'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.29%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.71'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '41.88'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.56%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.582.94'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -3.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.94'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.83%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '23.05'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.34%  '

# Restore default (unstyled) number format on all edited cells so their
# style attributes match the original (unstyled) cells.
$editedAddrs = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "E27", "D28", "E28", "B29", "C29", "D29", "E29", "B30", "C30", "D30", "E30", "D31", "E31", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "B39", "C39", "D39", "E39", "B40", "C40", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "E45", "B46", "C46", "D46", "E46", "B47", "C47", "D47", "E47", "D48", "E48", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $editedAddrs) {
    $fmtSource.Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

